$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9281321562220225
$ws.Range("D10").Value = 1.096363831534633
$ws.Range("E10").Value = 0.9768042872944713
$ws.Range("F10").Value = 1.014241451521715
$ws.Range("G10").Value = 0.9281321562220225
$ws.Range("H10").Value = 1.096363831534633
$ws.Range("I10").Value = 0.9635793264604577
$ws.Range("J10").Value = 1.014241451521715
$ws.Range("K10").Value = 0.9698176992801386
$ws.Range("L10").Value = 1.057161961425036
$ws.Range("M10").Value = 0.9281321562220225
$ws.Range("N10").Value = 1.036584059414552
$ws.Range("O10").Value = 1.00388543164321
$ws.Range("P10").Value = 1.002542770657524
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.8992523070443693
$ws.Range("D11").Value = 1.137092125381715
$ws.Range("E11").Value = 0.9639097732843847
$ws.Range("F11").Value = 1.024921279554959
$ws.Range("G11").Value = 0.8992523070443693
$ws.Range("H11").Value = 1.137092125381715
$ws.Range("I11").Value = 0.9446261236138204
$ws.Range("J11").Value = 1.026807574317993
$ws.Range("K11").Value = 0.9563097407852235
$ws.Range("L11").Value = 1.085497031555742
$ws.Range("M11").Value = 0.8992523070443693
$ws.Range("N11").Value = 1.05050094933305
$ws.Range("O11").Value = 1.006293871316357
$ws.Range("P11").Value = 1.004801994442276
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.8997916384301948
$ws.Range("D12").Value = 1.136299958894482
$ws.Range("E12").Value = 0.9641346366821307
$ws.Range("F12").Value = 1.024759567363152
$ws.Range("G12").Value = 0.8997916384301948
$ws.Range("H12").Value = 1.136299958894482
$ws.Range("I12").Value = 0.9448608516732403
$ws.Range("J12").Value = 1.026676126268792
$ws.Range("K12").Value = 0.9564495514641084
$ws.Range("L12").Value = 1.085056525735479
$ws.Range("M12").Value = 0.8997916384301948
$ws.Range("N12").Value = 1.050217297788306
$ws.Range("O12").Value = 1.00624645034249
$ws.Range("P12").Value = 1.004753607063947
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.8994271930749084
$ws.Range("D13").Value = 1.136914496889482
$ws.Range("E13").Value = 0.9639596816879765
$ws.Range("F13").Value = 1.024873319412216
$ws.Range("G13").Value = 0.8994271930749084
$ws.Range("H13").Value = 1.136914496889482
$ws.Range("I13").Value = 0.9446587769463682
$ws.Range("J13").Value = 1.026789059025635
$ws.Range("K13").Value = 0.9563291533961293
$ws.Range("L13").Value = 1.085428933178238
$ws.Range("M13").Value = 0.8994271930749084
$ws.Range("N13").Value = 1.050437089288729
$ws.Range("O13").Value = 1.006293672766146
$ws.Range("P13").Value = 1.004797576701369
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8126359999999988
$ws.Range("D14").Value = 1.264264
$ws.Range("E14").Value = 0.9262240000000009
$ws.Range("F14").Value = 1.053563999999999
$ws.Range("G14").Value = 0.8126359999999988
$ws.Range("H14").Value = 1.264264
$ws.Range("I14").Value = 0.9045520000000002
$ws.Range("J14").Value = 1.044599999999997
$ws.Range("K14").Value = 0.9328960000000011
$ws.Range("L14").Value = 1.154224
$ws.Range("M14").Value = 0.8126359999999988
$ws.Range("N14").Value = 1.095244
$ws.Range("O14").Value = 1.014172
$ws.Range("P14").Value = 1.01162
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.7
$ws.Range("D15").Value = 1.43
$ws.Range("E15").Value = 0.88
$ws.Range("F15").Value = 1.09
$ws.Range("G15").Value = 0.7
$ws.Range("H15").Value = 1.43
$ws.Range("I15").Value = 0.85
$ws.Range("J15").Value = 1.070687500000001
$ws.Range("K15").Value = 0.9
$ws.Range("L15").Value = 1.25
$ws.Range("M15").Value = 0.7
$ws.Range("N15").Value = 1.155
$ws.Range("O15").Value = 1.025
$ws.Range("P15").Value = 1.0213359375
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.8259820990464004
$ws.Range("D16").Value = 1.247448305459198
$ws.Range("E16").Value = 0.9279095173120047
$ws.Range("F16").Value = 1.0497095716864
$ws.Range("G16").Value = 0.8259820990464004
$ws.Range("H16").Value = 1.247448305459198
$ws.Range("I16").Value = 0.9101984391168084
$ws.Range("J16").Value = 1.039216147968002
$ws.Range("K16").Value = 0.9392661286912012
$ws.Range("L16").Value = 1.143353453977598
$ws.Range("M16").Value = 0.8259820990464004
$ws.Range("N16").Value = 1.087678911385601
$ws.Range("O16").Value = 1.012762373376001
$ws.Range("P16").Value = 1.010385457907201
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.994567703309683
$ws.Range("D17").Value = 0.9947379466564985
$ws.Range("E17").Value = 0.9956304470182653
$ws.Range("F17").Value = 0.9955278668947709
$ws.Range("G17").Value = 0.994567703309683
$ws.Range("H17").Value = 0.9947379466564985
$ws.Range("I17").Value = 0.994708758052603
$ws.Range("J17").Value = 0.9941770518350174
$ws.Range("K17").Value = 0.9943190310001944
$ws.Range("L17").Value = 0.9960287672887106
$ws.Range("M17").Value = 0.994567703309683
$ws.Range("N17").Value = 0.9951841968373819
$ws.Range("O17").Value = 0.9951159909698044
$ws.Range("P17").Value = 0.9949621965069679
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.004169490518929
$ws.Range("D18").Value = 0.9829809258565029
$ws.Range("E18").Value = 0.997742984701611
$ws.Range("F18").Value = 0.9927421262941738
$ws.Range("G18").Value = 1.004169490518929
$ws.Range("H18").Value = 0.9829809258565029
$ws.Range("I18").Value = 1.000530331442368
$ws.Range("J18").Value = 0.993145638405828
$ws.Range("K18").Value = 0.9984732363304983
$ws.Range("L18").Value = 0.9873534667360533
$ws.Range("M18").Value = 1.004169490518929
$ws.Range("N18").Value = 0.990361955279057
$ws.Range("O18").Value = 0.9944088818428041
$ws.Range("P18").Value = 0.9946422750357455
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.025918366506202
$ws.Range("D19").Value = 0.9506847712691185
$ws.Range("E19").Value = 1.005144779151771
$ws.Range("F19").Value = 0.9853445091929772
$ws.Range("G19").Value = 1.025918366506202
$ws.Range("H19").Value = 0.9506847712691185
$ws.Range("I19").Value = 1.011261434454208
$ws.Range("J19").Value = 0.986080596595171
$ws.Range("K19").Value = 1.00530962565929
$ws.Range("L19").Value = 0.965709638996423
$ws.Range("M19").Value = 1.025918366506202
$ws.Range("N19").Value = 0.9779147752104449
$ws.Range("O19").Value = 0.9917731065300173
$ws.Range("P19").Value = 0.9919317152281452

# Apply the existing HKL-index header style (bold, centered, thin border) to the
# newly added rows' A-column cells, matching the formatting already used by A2:A16.
$ws.Range("A2").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
